# "importação com separador funcionando"
# The sheet previously held numeric values in A1:A6 (a botched import where the
# field separator wasn't applied, so many numbers ended up in one column).
# Re-importing with the correct separator collapses the data down to three rows
# of text (each row keeps the original semicolon-separated values together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (4-6) from the old (wrong) import.
$ws.Range("A4:A6").EntireRow.Delete() | Out-Null

# Write the corrected, delimiter-aware rows as text values.
$ws.Range("A1").Value = "1520; 1595.6666; 2555"
$ws.Range("A2").Value = "1588,636; eyerte; 52525"
$ws.Range("A3").Value = "62626,5; 5858; 222"

# A3 used to carry a numeric style (custom number format); restore it to the
# default/normal style now that it holds text.
$ws.Range("A3").Style = "Normal"

# Move the active selection, matching where the user left off after the import.
$ws.Range("G8").Select() | Out-Null
